$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A288").Value = "IMX-USD"
$ws.Range("A289").Value = "TAO-USD"
$ws.Range("A290").Value = "MNT-USD"
